$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new data row for 2021-02-08 (serial 44235) ---
# This shifts the previously-existing rows 93-113 down to 94-114.
$ws.Rows.Item(93).Insert()

# Re-apply the column A-D formatting (borders / date number format) used by the
# surrounding data rows, since a freshly inserted row starts out unformatted.
$ws.Cells.Item(92,1).Copy()
$ws.Cells.Item(93,1).PasteSpecial(-4122)
$ws.Cells.Item(92,2).Copy()
$ws.Cells.Item(93,2).PasteSpecial(-4122)
$ws.Cells.Item(92,3).Copy()
$ws.Cells.Item(93,3).PasteSpecial(-4122)
$ws.Cells.Item(92,4).Copy()
$ws.Cells.Item(93,4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: write the refreshed date / nuovi-pos. / rolling-7gg values ---
# (rows 93-112; the 7-day rolling sum & rate in columns C/D shift along with
# the newly inserted day).
$ws.Cells.Item(93,1).Value = 44235
$ws.Cells.Item(93,2).Value = 0
$ws.Cells.Item(93,3).Value = 2
$ws.Cells.Item(93,4).Value = 96.15384615384616

$ws.Cells.Item(94,1).Value = 44236
$ws.Cells.Item(94,2).Value = 0
$ws.Cells.Item(94,3).Value = 1
$ws.Cells.Item(94,4).Value = 48.07692307692308

$ws.Cells.Item(95,1).Value = 44237
$ws.Cells.Item(95,2).Value = 0
$ws.Cells.Item(95,3).Value = 2
$ws.Cells.Item(95,4).Value = 96.15384615384616

$ws.Cells.Item(96,1).Value = 44238
$ws.Cells.Item(96,2).Value = 0
$ws.Cells.Item(96,3).Value = 2
$ws.Cells.Item(96,4).Value = 96.15384615384616

$ws.Cells.Item(97,1).Value = 44239
$ws.Cells.Item(97,2).Value = 0
$ws.Cells.Item(97,3).Value = 2
$ws.Cells.Item(97,4).Value = 96.15384615384616

$ws.Cells.Item(98,1).Value = 44240
$ws.Cells.Item(98,2).Value = 2
$ws.Cells.Item(98,3).Value = 2
$ws.Cells.Item(98,4).Value = 96.15384615384616

$ws.Cells.Item(99,1).Value = 44241
$ws.Cells.Item(99,2).Value = 0
$ws.Cells.Item(99,3).Value = 2
$ws.Cells.Item(99,4).Value = 96.15384615384616

$ws.Cells.Item(100,1).Value = 44242
$ws.Cells.Item(100,2).Value = 0
$ws.Cells.Item(100,3).Value = 2
$ws.Cells.Item(100,4).Value = 96.15384615384616

$ws.Cells.Item(101,1).Value = 44243
$ws.Cells.Item(101,2).Value = 0
$ws.Cells.Item(101,3).Value = 2
$ws.Cells.Item(101,4).Value = 96.15384615384616

$ws.Cells.Item(102,1).Value = 44244
$ws.Cells.Item(102,2).Value = 0
$ws.Cells.Item(102,3).Value = 1
$ws.Cells.Item(102,4).Value = 48.07692307692308

$ws.Cells.Item(103,1).Value = 44245
$ws.Cells.Item(103,2).Value = 0
$ws.Cells.Item(103,3).Value = 1
$ws.Cells.Item(103,4).Value = 48.07692307692308

$ws.Cells.Item(104,1).Value = 44246
$ws.Cells.Item(104,2).Value = 0
$ws.Cells.Item(104,3).Value = 1
$ws.Cells.Item(104,4).Value = 48.07692307692308

$ws.Cells.Item(105,1).Value = 44247
$ws.Cells.Item(105,2).Value = 1
$ws.Cells.Item(105,3).Value = 3
$ws.Cells.Item(105,4).Value = 144.2307692307692

$ws.Cells.Item(106,1).Value = 44248
$ws.Cells.Item(106,2).Value = 0
$ws.Cells.Item(106,3).Value = 3
$ws.Cells.Item(106,4).Value = 144.2307692307692

$ws.Cells.Item(107,1).Value = 44249
$ws.Cells.Item(107,2).Value = 0
$ws.Cells.Item(107,3).Value = 4
$ws.Cells.Item(107,4).Value = 192.3076923076923

$ws.Cells.Item(108,1).Value = 44250
$ws.Cells.Item(108,2).Value = 2
$ws.Cells.Item(108,3).Value = 4
$ws.Cells.Item(108,4).Value = 192.3076923076923

$ws.Cells.Item(109,1).Value = 44251
$ws.Cells.Item(109,2).Value = 0
$ws.Cells.Item(109,3).Value = 5
$ws.Cells.Item(109,4).Value = 240.3846153846154

$ws.Cells.Item(110,1).Value = 44252
$ws.Cells.Item(110,2).Value = 1
$ws.Cells.Item(110,3).Value = 5
$ws.Cells.Item(110,4).Value = 240.3846153846154

$ws.Cells.Item(111,1).Value = 44253
$ws.Cells.Item(111,2).Value = 0
$ws.Cells.Item(111,3).Value = 6
$ws.Cells.Item(111,4).Value = 288.4615384615385

$ws.Cells.Item(112,1).Value = 44254
$ws.Cells.Item(112,2).Value = 2
$ws.Cells.Item(112,3).Value = 5
$ws.Cells.Item(112,4).Value = 240.3846153846154

# Rows 113 and 114 (2021-02-28 / 2021-03-01, serials 44255/44256) already carry the
# correct values after the insert above, so they are left untouched.

# --- Step 3: append a brand-new final row for 2021-03-02 (serial 44257) ---
$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 1

# Match the formatting of the preceding row for the new trailing row.
$ws.Cells.Item(114,1).Copy()
$ws.Cells.Item(115,1).PasteSpecial(-4122)
$ws.Cells.Item(114,3).Copy()
$ws.Cells.Item(115,3).PasteSpecial(-4122)
$ws.Cells.Item(114,4).Copy()
$ws.Cells.Item(115,4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-assert the date/value after the format paste (PasteSpecial only touches formatting).
$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 1
